$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "M1"
$ws.Cells.Item(2, 2).Value = "ECs"
$ws.Cells.Item(2, 3).Value = "F3"
$ws.Cells.Item(2, 4).Value = "M1"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 1.341690333333333
$ws.Cells.Item(2, 8).Value = 4.025071000000001
$ws.Cells.Item(2, 9).Value = 0.05354309921306734
$ws.Cells.Item(2, 10).Value = 0.05354309921306735
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.9880003333333333
$ws.Cells.Item(2, 14).Value = 2.964001
$ws.Cells.Item(2, 15).Value = 0.02006000579181712
$ws.Cells.Item(2, 16).Value = 0.02006000579181712
$ws.Cells.Item(2, 17).Value = 1.325590496563444
$ws.Cells.Item(2, 18).Value = 11.930314469071
$ws.Cells.Item(2, 19).Value = 0.00107407488032597
$ws.Cells.Item(2, 20).Value = 0.00107407488032597

# Row 3
$ws.Cells.Item(3, 1).Value = "M1"
$ws.Cells.Item(3, 2).Value = "ECs"
$ws.Cells.Item(3, 3).Value = "F3"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 1.341690333333333
$ws.Cells.Item(3, 8).Value = 4.025071000000001
$ws.Cells.Item(3, 9).Value = 0.05354309921306734
$ws.Cells.Item(3, 10).Value = 0.05354309921306735
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 46.102515
$ws.Cells.Item(3, 14).Value = 138.307545
$ws.Cells.Item(3, 15).Value = 0.9360489938269279
$ws.Cells.Item(3, 16).Value = 0.9360489938269277
$ws.Cells.Item(3, 17).Value = 61.85529871785501
$ws.Cells.Item(3, 18).Value = 556.6976884606951
$ws.Cells.Item(3, 19).Value = 0.05011896414476705
$ws.Cells.Item(3, 20).Value = 0.05011896414476706

# Row 4
$ws.Cells.Item(4, 1).Value = "M1"
$ws.Cells.Item(4, 2).Value = "ECs"
$ws.Cells.Item(4, 3).Value = "F3"
$ws.Cells.Item(4, 4).Value = "M2"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 1.341690333333333
$ws.Cells.Item(4, 8).Value = 4.025071000000001
$ws.Cells.Item(4, 9).Value = 0.05354309921306734
$ws.Cells.Item(4, 10).Value = 0.05354309921306735
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 0.6702083333333334
$ws.Cells.Item(4, 14).Value = 2.010625
$ws.Cells.Item(4, 15).Value = 0.01360767055921112
$ws.Cells.Item(4, 16).Value = 0.01360767055921112
$ws.Cells.Item(4, 17).Value = 0.899212042152778
$ws.Cells.Item(4, 18).Value = 8.092908379375002
$ws.Cells.Item(4, 19).Value = 0.0007285968548105765
$ws.Cells.Item(4, 20).Value = 0.0007285968548105766

# Row 5
$ws.Cells.Item(5, 1).Value = "M1"
$ws.Cells.Item(5, 2).Value = "ECs"
$ws.Cells.Item(5, 3).Value = "F3"
$ws.Cells.Item(5, 4).Value = "F7"
$ws.Cells.Item(5, 5).Value = 2
$ws.Cells.Item(5, 6).Value = 0.6666666666666666
$ws.Cells.Item(5, 7).Value = 1.341690333333333
$ws.Cells.Item(5, 8).Value = 4.025071000000001
$ws.Cells.Item(5, 9).Value = 0.05354309921306734
$ws.Cells.Item(5, 10).Value = 0.05354309921306735
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.794312
$ws.Cells.Item(5, 14).Value = 2.382936
$ws.Cells.Item(5, 15).Value = 0.01612742706953525
$ws.Cells.Item(5, 16).Value = 0.01612742706953524
$ws.Cells.Item(5, 17).Value = 1.065720732050667
$ws.Cells.Item(5, 18).Value = 9.591486588456
$ws.Cells.Item(5, 19).Value = 0.0008635124276356335
$ws.Cells.Item(5, 20).Value = 0.0008635124276356335

# Row 6
$ws.Cells.Item(6, 1).Value = "M1"
$ws.Cells.Item(6, 2).Value = "ECs"
$ws.Cells.Item(6, 3).Value = "F3"
$ws.Cells.Item(6, 4).Value = "sCs"
$ws.Cells.Item(6, 5).Value = 2
$ws.Cells.Item(6, 6).Value = 0.6666666666666666
$ws.Cells.Item(6, 7).Value = 1.341690333333333
$ws.Cells.Item(6, 8).Value = 4.025071000000001
$ws.Cells.Item(6, 9).Value = 0.05354309921306734
$ws.Cells.Item(6, 10).Value = 0.05354309921306735
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 0.69721
$ws.Cells.Item(6, 14).Value = 2.09163
$ws.Cells.Item(6, 15).Value = 0.01415590275250867
$ws.Cells.Item(6, 16).Value = 0.01415590275250867
$ws.Cells.Item(6, 17).Value = 0.9354399173033334
$ws.Cells.Item(6, 18).Value = 8.41895925573
$ws.Cells.Item(6, 19).Value = 0.0007579509055281048
$ws.Cells.Item(6, 20).Value = 0.0007579509055281049

# Row 7
$ws.Cells.Item(7, 1).Value = "M2"
$ws.Cells.Item(7, 2).Value = "ECs"
$ws.Cells.Item(7, 3).Value = "F3"
$ws.Cells.Item(7, 4).Value = "M1"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 15.617696
$ws.Cells.Item(7, 8).Value = 46.853088
$ws.Cells.Item(7, 9).Value = 0.6232584566142001
$ws.Cells.Item(7, 10).Value = 0.6232584566142001
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 0.9880003333333333
$ws.Cells.Item(7, 14).Value = 2.964001
$ws.Cells.Item(7, 15).Value = 0.02006000579181712
$ws.Cells.Item(7, 16).Value = 0.02006000579181712
$ws.Cells.Item(7, 17).Value = 15.43028885389867
$ws.Cells.Item(7, 18).Value = 138.872599685088
$ws.Cells.Item(7, 19).Value = 0.01250256824947986
$ws.Cells.Item(7, 20).Value = 0.01250256824947985

# Row 8
$ws.Cells.Item(8, 1).Value = "M2"
$ws.Cells.Item(8, 2).Value = "ECs"
$ws.Cells.Item(8, 3).Value = "F3"
$ws.Cells.Item(8, 4).Value = "FAPs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 15.617696
$ws.Cells.Item(8, 8).Value = 46.853088
$ws.Cells.Item(8, 9).Value = 0.6232584566142001
$ws.Cells.Item(8, 10).Value = 0.6232584566142001
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 46.102515
$ws.Cells.Item(8, 14).Value = 138.307545
$ws.Cells.Item(8, 15).Value = 0.9360489938269279
$ws.Cells.Item(8, 16).Value = 0.9360489938269277
$ws.Cells.Item(8, 17).Value = 720.01506410544
$ws.Cells.Item(8, 18).Value = 6480.13557694896
$ws.Cells.Item(8, 19).Value = 0.583400451207846
$ws.Cells.Item(8, 20).Value = 0.5834004512078459

# Row 9
$ws.Cells.Item(9, 1).Value = "M2"
$ws.Cells.Item(9, 2).Value = "ECs"
$ws.Cells.Item(9, 3).Value = "F3"
$ws.Cells.Item(9, 4).Value = "M2"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 15.617696
$ws.Cells.Item(9, 8).Value = 46.853088
$ws.Cells.Item(9, 9).Value = 0.6232584566142001
$ws.Cells.Item(9, 10).Value = 0.6232584566142001
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 0.6702083333333334
$ws.Cells.Item(9, 14).Value = 2.010625
$ws.Cells.Item(9, 15).Value = 0.01360767055921112
$ws.Cells.Item(9, 16).Value = 0.01360767055921112
$ws.Cells.Item(9, 17).Value = 10.46711000666667
$ws.Cells.Item(9, 18).Value = 94.20399006000001
$ws.Cells.Item(9, 19).Value = 0.00848109575084841
$ws.Cells.Item(9, 20).Value = 0.00848109575084841

# Row 10
$ws.Cells.Item(10, 1).Value = "M2"
$ws.Cells.Item(10, 2).Value = "ECs"
$ws.Cells.Item(10, 3).Value = "F3"
$ws.Cells.Item(10, 4).Value = "F7"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 15.617696
$ws.Cells.Item(10, 8).Value = 46.853088
$ws.Cells.Item(10, 9).Value = 0.6232584566142001
$ws.Cells.Item(10, 10).Value = 0.6232584566142001
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 0.794312
$ws.Cells.Item(10, 14).Value = 2.382936
$ws.Cells.Item(10, 15).Value = 0.01612742706953525
$ws.Cells.Item(10, 16).Value = 0.01612742706953524
$ws.Cells.Item(10, 17).Value = 12.405323345152
$ws.Cells.Item(10, 18).Value = 111.647910106368
$ws.Cells.Item(10, 19).Value = 0.01005155530451661
$ws.Cells.Item(10, 20).Value = 0.01005155530451661

# Row 11
$ws.Cells.Item(11, 1).Value = "M2"
$ws.Cells.Item(11, 2).Value = "ECs"
$ws.Cells.Item(11, 3).Value = "F3"
$ws.Cells.Item(11, 4).Value = "sCs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 15.617696
$ws.Cells.Item(11, 8).Value = 46.853088
$ws.Cells.Item(11, 9).Value = 0.6232584566142001
$ws.Cells.Item(11, 10).Value = 0.6232584566142001
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 0.69721
$ws.Cells.Item(11, 14).Value = 2.09163
$ws.Cells.Item(11, 15).Value = 0.01415590275250867
$ws.Cells.Item(11, 16).Value = 0.01415590275250867
$ws.Cells.Item(11, 17).Value = 10.88881382816
$ws.Cells.Item(11, 18).Value = 97.99932445344
$ws.Cells.Item(11, 19).Value = 0.008822786101509262
$ws.Cells.Item(11, 20).Value = 0.00882278610150926

# Row 12
$ws.Cells.Item(12, 1).Value = "F7"
$ws.Cells.Item(12, 2).Value = "ECs"
$ws.Cells.Item(12, 3).Value = "F3"
$ws.Cells.Item(12, 4).Value = "M1"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 8.098751
$ws.Cells.Item(12, 8).Value = 24.296253
$ws.Cells.Item(12, 9).Value = 0.3231984441727326
$ws.Cells.Item(12, 10).Value = 0.3231984441727326
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 0.9880003333333333
$ws.Cells.Item(12, 14).Value = 2.964001
$ws.Cells.Item(12, 15).Value = 0.02006000579181712
$ws.Cells.Item(12, 16).Value = 0.02006000579181712
$ws.Cells.Item(12, 17).Value = 8.001568687583665
$ws.Cells.Item(12, 18).Value = 72.014118188253
$ws.Cells.Item(12, 19).Value = 0.006483362662011299
$ws.Cells.Item(12, 20).Value = 0.006483362662011299

# Row 13
$ws.Cells.Item(13, 1).Value = "F7"
$ws.Cells.Item(13, 2).Value = "ECs"
$ws.Cells.Item(13, 3).Value = "F3"
$ws.Cells.Item(13, 4).Value = "FAPs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 8.098751
$ws.Cells.Item(13, 8).Value = 24.296253
$ws.Cells.Item(13, 9).Value = 0.3231984441727326
$ws.Cells.Item(13, 10).Value = 0.3231984441727326
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 46.102515
$ws.Cells.Item(13, 14).Value = 138.307545
$ws.Cells.Item(13, 15).Value = 0.9360489938269279
$ws.Cells.Item(13, 16).Value = 0.9360489938269277
$ws.Cells.Item(13, 17).Value = 373.3727894587651
$ws.Cells.Item(13, 18).Value = 3360.355105128885
$ws.Cells.Item(13, 19).Value = 0.3025295784743148
$ws.Cells.Item(13, 20).Value = 0.3025295784743149

# Row 14
$ws.Cells.Item(14, 1).Value = "F7"
$ws.Cells.Item(14, 2).Value = "ECs"
$ws.Cells.Item(14, 3).Value = "F3"
$ws.Cells.Item(14, 4).Value = "M2"
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 8.098751
$ws.Cells.Item(14, 8).Value = 24.296253
$ws.Cells.Item(14, 9).Value = 0.3231984441727326
$ws.Cells.Item(14, 10).Value = 0.3231984441727326
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 0.6702083333333334
$ws.Cells.Item(14, 14).Value = 2.010625
$ws.Cells.Item(14, 15).Value = 0.01360767055921112
$ws.Cells.Item(14, 16).Value = 0.01360767055921112
$ws.Cells.Item(14, 17).Value = 5.427850409791668
$ws.Cells.Item(14, 18).Value = 48.850653688125
$ws.Cells.Item(14, 19).Value = 0.004397977953552131
$ws.Cells.Item(14, 20).Value = 0.004397977953552131

# Row 15
$ws.Cells.Item(15, 1).Value = "F7"
$ws.Cells.Item(15, 2).Value = "ECs"
$ws.Cells.Item(15, 3).Value = "F3"
$ws.Cells.Item(15, 4).Value = "F7"
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 8.098751
$ws.Cells.Item(15, 8).Value = 24.296253
$ws.Cells.Item(15, 9).Value = 0.3231984441727326
$ws.Cells.Item(15, 10).Value = 0.3231984441727326
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 0.794312
$ws.Cells.Item(15, 14).Value = 2.382936
$ws.Cells.Item(15, 15).Value = 0.01612742706953525
$ws.Cells.Item(15, 16).Value = 0.01612742706953524
$ws.Cells.Item(15, 17).Value = 6.432935104312
$ws.Cells.Item(15, 18).Value = 57.896415938808
$ws.Cells.Item(15, 19).Value = 0.005212359337383003
$ws.Cells.Item(15, 20).Value = 0.005212359337383003

# Row 16
$ws.Cells.Item(16, 1).Value = "F7"
$ws.Cells.Item(16, 2).Value = "ECs"
$ws.Cells.Item(16, 3).Value = "F3"
$ws.Cells.Item(16, 4).Value = "sCs"
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 8.098751
$ws.Cells.Item(16, 8).Value = 24.296253
$ws.Cells.Item(16, 9).Value = 0.3231984441727326
$ws.Cells.Item(16, 10).Value = 0.3231984441727326
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 0.69721
$ws.Cells.Item(16, 14).Value = 2.09163
$ws.Cells.Item(16, 15).Value = 0.01415590275250867
$ws.Cells.Item(16, 16).Value = 0.01415590275250867
$ws.Cells.Item(16, 17).Value = 5.64653018471
$ws.Cells.Item(16, 18).Value = 50.81877166239
$ws.Cells.Item(16, 19).Value = 0.004575165745471305
$ws.Cells.Item(16, 20).Value = 0.004575165745471305
